$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "23.276.59"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.16%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.621.57"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("E4").Value = "  +0.75%  "

$ws.Range("E5").Value = "  +0.58%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "303.90"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.32%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3795"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("E8").Value = "  -2.33%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3602"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.08%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.222"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -5.45%  "

$ws.Range("E11").Value = "  -1.73%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "22.54"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.70%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.531"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.41%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.00001241"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.64%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "7.207"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.69%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.621.70"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.54%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "93.39"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.65%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06904"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.78%  "

$ws.Range("E20").Value = "  -3.38%  "

$ws.Range("E21").Value = "  +0.59%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.407"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.00%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "23.288.45"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.14%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "12.69"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.62%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.176"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.451"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "21.06"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.84%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "149.38"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.16%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.280"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.73%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "134.36"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.74%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.289"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.48%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.802.90"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "6.745"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.65%  "

$ws.Range("E34").Value = "  +3.83%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.9448"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.11%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02789"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.2511"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.32%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.08819"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.58%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "6.032"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.12%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.07096"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -5.28%  "

$ws.Range("E41").Value = "  -3.81%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.7006"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.65%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "16.06"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.90%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "12.22"
$c.Style = "Normal"

$ws.Range("E45").Value = "  +0.47%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.6409"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.72%  "

$ws.Range("E47").Value = "  -2.81%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.984"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.37%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.07965"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.12%  "

$ws.Range("E50").Value = "  -2.05%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "125.07"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -5.92%  "
